# Horarios actualizados Línea 141 - 674
$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912" (sheet1): new scrape with 3 data rows ----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:57:27"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6 - updated in place
$ws1.Cells.Item(6, 1).Value = "01:57:27"
$ws1.Cells.Item(6, 2).Value = "02:00"
$ws1.Cells.Item(6, 3).Value = "14_ABASTO"
$ws1.Cells.Item(6, 4).Value = 3
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Row 7 - new
$ws1.Cells.Item(7, 1).Value = "01:57:27"
$ws1.Cells.Item(7, 2).Value = "02:58"
$ws1.Cells.Item(7, 3).Value = "215_ALUAR"
$ws1.Cells.Item(7, 4).Value = 61
$ws1.Cells.Item(7, 5).Value = "LP1912"

# Row 8 - new
$ws1.Cells.Item(8, 1).Value = "01:57:27"
$ws1.Cells.Item(8, 2).Value = "03:50"
$ws1.Cells.Item(8, 3).Value = "14_ABASTO"
$ws1.Cells.Item(8, 4).Value = 113
$ws1.Cells.Item(8, 5).Value = "LP1912"

# ---- Sheet "LP1912-215" (sheet2): timestamp + row 6 refresh ----
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:57:27"

$ws2.Cells.Item(6, 1).Value = "01:57:27"
$ws2.Cells.Item(6, 2).Value = "02:58"
$ws2.Cells.Item(6, 4).Value = 61

# ---- Sheet "6203-6173" (sheet3): timestamp refresh only ----
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:57:27"
